$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet currently holds an 11-column x 9-row table (A1:K9): a header row
# of car names across the top and color names down the left column, with
# percentages in the body.
#
# This transposes the table in place into a 9-column x 11-row table
# (A1:I11): color names now run across the header row and car names run
# down the left column. The top-left header label also changes from
# "Car Name" to "Sheet Name".
# ---------------------------------------------------------------------------

$oldRows = 9
$oldCols = 11
$stageRowOffset = 1000

# 1) Stage a copy of the original grid well out of the way so that writing
#    the transposed values back over the original A1:K9 block can never
#    clobber a value we still need to read (the two ranges overlap in the
#    top-left 9x9 square).
for ($r = 1; $r -le $oldRows; $r++) {
    for ($c = 1; $c -le $oldCols; $c++) {
        $val = $ws.Cells.Item($r, $c).Value()
        $ws.Cells.Item($stageRowOffset + $r, $c).Value = $val
    }
}

# 2) Clear the entire original block - the new table has a different shape
#    (9 cols x 11 rows instead of 11 cols x 9 rows) so the two trailing
#    columns (J:K) fall outside the new dimension and must not survive.
$ws.Range("A1:K9").Clear()

# 3) Write the transposed values back: new row i, col j <- old row j, col i.
for ($r = 1; $r -le $oldRows; $r++) {
    for ($c = 1; $c -le $oldCols; $c++) {
        $val = $ws.Cells.Item($stageRowOffset + $r, $c).Value()
        $ws.Cells.Item($c, $r).Value = $val
    }
}

# 4) Remove the staged scratch copy.
$ws.Range("A1001:K1009").Clear()

# 5) The header of the new left-hand column is "Sheet Name", not the old
#    "Car Name" label that the transpose carried over from A1.
$ws.Range("A1").Value = "Sheet Name"
